$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the oldest 10 years of data (2000年-2009年), which are currently
# in rows 2 through 11. Deleting these rows shifts rows 12-22
# (2010年-2020年) up to become rows 2-12.
$ws.Rows("2:11").Delete()

# Append the new year of data (2021年) as the new last row (row 13).
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 78024

# Match the existing formatting used for the other year rows (border,
# bold, centered alignment) by copying the format from the row above.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

# Update the sheet dimension/selection to reflect the new data extent.
$ws.Range("A1:B13").Select()
